$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("0", 4.5, 31.251),
    @("1", 5.1, 31.389),
    @("2", 5, 31.739),
    @("3", 3.9, 31.404),
    @("4", 6.1, 31.262),
    @("5", 5.6, 31.379),
    @("7", 4.5, 31.593),
    @("9", 4.7, 31.295),
    @("11", 4.7, 31.627),
    @("13", 4, 31.494),
    @("15", 4.8, 31.469),
    @("17", 4.5, 31.333),
    @("19", 3.7, 31.324),
    @("21", 4.4, 31.251),
    @("23", 3.9, 31.904),
    @("25", 4.4, 31.621),
    @("27", 6, 31.383),
    @("29", 5, 31.52),
    @("31", 5.4, 31.502),
    @("33", 3.7, 31.636),
    @("35", 4.7, 31.51),
    @("37", 5.4, 31.853),
    @("39", 3.9, 31.441),
    @("41", 4.8, 31.691),
    @("43", 3.9, 31.566),
    @("45", 4.4, 31.473)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $vals = $data[$i]
    $ws.Cells.Item($row, 1).Value = $vals[0]
    $ws.Cells.Item($row, 2).Value = $vals[1]
    $ws.Cells.Item($row, 3).Value = $vals[2]
}
